$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.258.30"
$ws.Range("E2").Value = "  +2.08%  "

# Row 3
$ws.Range("D3").Value = "2.349.03"
$ws.Range("E3").Value = "  +6.08%  "

# Row 4
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
$ws.Range("D5").Value = "314.50"
$ws.Range("E5").Value = "  +6.43%  "

# Row 6
$ws.Range("D6").Value = "109.54"
$ws.Range("E6").Value = "  +1.94%  "

# Row 7
$ws.Range("D7").Value = "0.643"
$ws.Range("E7").Value = "  +3.33%  "

# Row 8
$ws.Range("E8").Value = "  -0.23%  "

# Row 9
$ws.Range("D9").Value = "0.634"
$ws.Range("E9").Value = "  +6.63%  "

# Row 10
$ws.Range("D10").Value = "42.98"
$ws.Range("E10").Value = "  -1.16%  "

# Row 11
$ws.Range("D11").Value = "0.0939"
$ws.Range("E11").Value = "  +3.31%  "

# Row 12
$ws.Range("E12").Value = "  +1.48%  "

# Row 13
$ws.Range("D13").Value = "1.04"
$ws.Range("E13").Value = "  +9.19%  "

# Row 14
$ws.Range("E14").Value = "  +2.29%  "

# Row 15
$ws.Range("D15").Value = "16.28"
$ws.Range("E15").Value = "  +9.26%  "

# Row 16
$ws.Range("D16").Value = "2.705.73"
$ws.Range("E16").Value = "  +6.19%  "

# Row 17
$ws.Range("D17").Value = "2.343.51"
$ws.Range("E17").Value = "  +4.90%  "

# Row 18
$ws.Range("D18").Value = "43.260.52"
$ws.Range("E18").Value = "  +2.25%  "

# Row 19
$ws.Range("E19").Value = "  +3.40%  "

# Row 20
$ws.Range("D20").Value = "7.26"
$ws.Range("E20").Value = "  -1.36%  "

# Row 21
$ws.Range("D21").Value = "75.48"
$ws.Range("E21").Value = "  +3.98%  "

# Row 22
$ws.Range("B22").Value = "PancakeSwap"
$ws.Range("C22").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D22").Value = "3.46"
$ws.Range("E22").Value = "  +0.36%  "

# Row 23
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").Value = "2.58"
$ws.Range("E23").Value = "  +13.06%  "

# Row 24
$ws.Range("D24").Value = "254.99"
$ws.Range("E24").Value = "  +12.10%  "

# Row 25
$ws.Range("D25").Value = "9.13"
$ws.Range("E25").Value = "  +1.23%  "

# Row 26
$ws.Range("D26").Value = "12.05"
$ws.Range("E26").Value = "  +4.37%  "

# Row 27
$ws.Range("E27").Value = "  +0.03%  "

# Row 28
$ws.Range("D28").Value = "39.21"
$ws.Range("E28").Value = "  +1.89%  "

# Row 29
$ws.Range("E29").Value = "  +1.65%  "

# Row 30
$ws.Range("D30").Value = "22.37"

# Row 31
$ws.Range("D31").Value = "174.88"
$ws.Range("E31").Value = "  +0.90%  "

# Row 32
$ws.Range("E32").Value = "  -1.13%  "

# Row 33
$ws.Range("E33").Value = "  +4.92%  "

# Row 34
$ws.Range("D34").Value = "6.05"
$ws.Range("E34").Value = "  +9.99%  "

# Row 35
$ws.Range("D35").Value = "0.133"
$ws.Range("E35").Value = "  +6.19%  "

# Row 36
$ws.Range("E36").Value = "  -1.15%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.0378"
$ws.Range("E37").Value = "  +4.35%  "

# Row 38
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "4.13"
$ws.Range("E38").Value = "  -3.77%  "

# Row 39
$ws.Range("E39").Value = "  +2.17%  "

# Row 40
$ws.Range("E40").Value = "  +11.35%  "

# Row 41
$ws.Range("D41").Value = "72.81"
$ws.Range("E41").Value = "  +2.22%  "

# Row 42
$ws.Range("E42").Value = "  +14.66%  "

# Row 43
$ws.Range("D43").Value = "0.234"
$ws.Range("E43").Value = "  +2.02%  "

# Row 44
$ws.Range("D44").Value = "12.80"
$ws.Range("E44").Value = "  +1.72%  "

# Row 45
$ws.Range("E45").Value = "  -0.08%  "

# Row 46
$ws.Range("D46").Value = "5.63"
$ws.Range("E46").Value = "  +3.97%  "

# Row 47
$ws.Range("D47").Value = "9.32"
$ws.Range("E47").Value = "  +11.76%  "

# Row 48
$ws.Range("D48").Value = "111.07"
$ws.Range("E48").Value = "  +7.79%  "

# Row 49
$ws.Range("E49").Value = "  -0.05%  "

# Row 50
$ws.Range("E50").Value = "  +3.94%  "

# Row 51
$ws.Range("D51").Value = "69.84"
$ws.Range("E51").Value = "  +5.72%  "
